$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 = 3)
$ws.Range("B2").Value = 0.34
$ws.Range("C2").Value = 1.01
$ws.Range("D2").Value = 0.4
$ws.Range("E2").Value = 0.03
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0

# Row 3 (A3 = 4)
$ws.Range("B3").Value = 0.34
$ws.Range("C3").Value = 1.01
$ws.Range("D3").Value = 0.4
$ws.Range("E3").Value = 0.29
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.01

# Row 4 (A4 = 5): update B4, C4, D4; remove E4, F4, G4 entirely
$ws.Range("B4").Value = 0.34
$ws.Range("C4").Value = 1.01
$ws.Range("D4").Value = 0.4
$ws.Range("E4:G4").ClearContents()
